# Apply the daily crypto-price/volume refresh to Sheet1.
# Most D/E cells are plain text (prices with thousands-dot separators and
# padded percentage strings), but some new D-column values parse as plain
# decimals (e.g. "328.26"). Assigning those straight to .Value would let
# Excel auto-coerce the text to a number (losing formatting / trailing
# zeros and introducing float noise), so those cells are briefly switched
# to text format, written, then had their format cleared again so the
# cell keeps its original (default) style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.140.11'
$ws.Range("E2").Value = '  -0.34%  '
$ws.Range("D3").Value = '2.403.96'
$ws.Range("E3").Value = '  +5.63%  '
$ws.Range("E4").Value = '  -0.31%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '328.26'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +8.41%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '105.88'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -7.35%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.655'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +3.00%  '
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.654'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +6.15%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.20'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -5.72%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0941'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.45%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.78'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.32%  '
$ws.Range("E13").Value = '  +0.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '17.15'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +11.03%  '
$ws.Range("E15").Value = '  +1.96%  '
$ws.Range("D16").Value = '2.770.36'
$ws.Range("E16").Value = '  +5.82%  '
$ws.Range("D17").Value = '2.396.63'
$ws.Range("E17").Value = '  +5.18%  '
$ws.Range("D18").Value = '43.282.55'
$ws.Range("E18").Value = '  +0.18%  '
$ws.Range("E19").Value = '  +7.39%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000108'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '77.18'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.50%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.77'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +6.32%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '273.98'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +6.73%  '
$ws.Range("E24").Value = '  -0.41%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.75'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +7.22%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.99'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +2.25%  '
$ws.Range("E27").Value = '  +0.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.18'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +4.14%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '176.64'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.71%  '
$ws.Range("E30").Value = '  -2.21%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '37.27'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -2.73%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0942'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +4.70%  '
$ws.Range("E33").Value = '  -0.71%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.98'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +4.22%  '
$ws.Range("E35").Value = '  +5.19%  '
$ws.Range("E36").Value = '  -3.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.14'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -3.42%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0365'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -3.36%  '
$ws.Range("E39").Value = '  +2.62%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.86'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +15.16%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.59'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +17.00%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.236'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +1.17%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '70.37'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.85%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '123.50'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +14.94%  '
$ws.Range("E45").Value = '  +0.02%  '
$ws.Range("B46").Value = 'BitcoinSV'
$ws.Range("C46").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '91.31'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +43.69%  '
$ws.Range("B47").Value = 'Celestia'
$ws.Range("C47").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.36'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.76%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.60'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.73%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.36'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +6.70%  '
$ws.Range("E50").Value = '  +1.68%  '
$ws.Range("E51").Value = '  +12.63%  '
